# QA Compiler: Master_Knowledge STATUS sheet - transfer mock QA results in.
#
# The "Bob" row (row 2) is replaced by the results that used to live on the
# "Doni" row (row 3), and the now-duplicate row 3 is deleted, shrinking the
# sheet from A1:F3 down to A1:F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STATUS")

$xlPasteFormats = -4122

# Columns C/D are plain numbers - safe to assign directly.
$ws.Cells.Item(2, 1).Value = "Doni"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0

# Columns B/E/F store percentages as literal text (e.g. "0%"), not numbers.
# Assigning that string straight to .Value makes Excel parse it as a real
# percentage, which would also rewrite the cell's number format. Flip the
# cell to text format first so the literal string sticks, then paste the
# original formatting back from column C (untouched, same row) so the
# cell's look (border/fill/number format) ends up exactly as it was.
foreach ($col in 2, 5, 6) {
    $cell = $ws.Cells.Item(2, $col)
    $cell.NumberFormat = "@"
    $cell.Value = "0%"
    $ws.Cells.Item(2, 3).Copy()
    $cell.PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# Row 3 ("Doni") has now been folded into row 2, so drop the old row.
$ws.Rows.Item(3).Delete()
